# Add an "Appear" entrance animation (triggered on click) to the table
# shape ("Tabelle 4", shape id 5) on slide 3 ("Tests und Bug-Fixes").
#
# msoAnimEffectAppear = 1 (effectId)
# msoAnimTriggerOnPageClick = 1 (default trigger, explicit here for clarity)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item("Tabelle 4")

$seq = $s.TimeLine.MainSequence
$effect = $seq.AddEffect($sh, 1, 0, 1)
